# Add a "Sanity" execution sheet (a cut-down subset of the existing
# "Scenarios" sheet) so the suite's execution sheet can be switched at
# run time between a quick sanity pass and the full regression pass.

$wb = $excel.ActiveWorkbook

# Duplicate the "Scenarios" sheet, inserting the copy right before it.
# This becomes the new "Sanity" sheet; the original is kept (and later
# renamed) as the "Regression" sheet.
$scenarios = $wb.Worksheets.Item("Scenarios")
$scenarios.Copy($scenarios)

# Re-resolve fresh references to the sheets by name: after Copy() the
# old `$scenarios` handle tracks the tab *position*, not the original
# sheet object, so it now points at the freshly inserted copy.
$sanity = $wb.Worksheets.Item("Scenarios (2)")
$regression = $wb.Worksheets.Item("Scenarios")

$sanity.Name = "Sanity"
$regression.Name = "Regression"

# On the Sanity sheet, only keep the AccountLoginTest scenario (row 3)
# enabled; switch every other scenario's Run_Mode to "No" so only a
# minimal smoke test executes when this sheet is selected.
$sanity.Range("D2").Value = "No"
$sanity.Range("D4").Value = "No"
$sanity.Range("D5").Value = "No"
$sanity.Range("D6").Value = "No"
$sanity.Range("D7").Value = "No"
$sanity.Range("D8").Value = "No"
$sanity.Range("D9").Value = "No"
$sanity.Range("D10").Value = "No"
$sanity.Range("D11").Value = "No"

# Leave the selection/active cell on each sheet the way the author did,
# with Sanity as the active (first, selected) tab.
$regression.Activate() | Out-Null
$regression.Range("E17").Select() | Out-Null

$sanity.Activate() | Out-Null
$sanity.Range("C13").Select() | Out-Null
